# Updates the cryptocurrency price/volume table with refreshed market data.
# Numeric-looking Price values are written with a leading apostrophe so
# Excel stores them as text (matching the sheet's original text-cell
# formatting for the "Price" column) instead of auto-converting them to
# real numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '95.632.96'
$ws.Range("E2").Value = '  -1.52%  '

# Row 3
$ws.Range("D3").Value = '3.607.35'
$ws.Range("E3").Value = '  -2.61%  '

# Row 4
$ws.Range("D4").Value = '''2.73'
$ws.Range("E4").Value = '  +32.94%  '

# Row 5
$ws.Range("D5").Value = '''0.998'
$ws.Range("E5").Value = '  -0.22%  '

# Row 6
$ws.Range("D6").Value = '''222.56'
$ws.Range("E6").Value = '  -5.63%  '

# Row 7
$ws.Range("D7").Value = '''638.11'
$ws.Range("E7").Value = '  -2.81%  '

# Row 8
$ws.Range("D8").Value = '''0.424'
$ws.Range("E8").Value = '  -1.64%  '

# Row 9
$ws.Range("D9").Value = '''1.19'
$ws.Range("E9").Value = '  +10.17%  '

# Row 10
$ws.Range("E10").Value = '  -0.06%  '

# Row 11
$ws.Range("D11").Value = '3.606.27'
$ws.Range("E11").Value = '  -2.58%  '

# Row 12
$ws.Range("D12").Value = '''48.28'
$ws.Range("E12").Value = '  +7.06%  '

# Row 13
$ws.Range("E13").Value = '  +3.02%  '

# Row 14
$ws.Range("D14").Value = '''0.0000293'
$ws.Range("E14").Value = '  -5.33%  '

# Row 15
$ws.Range("D15").Value = '''6.51'
$ws.Range("E15").Value = '  -5.18%  '

# Row 16
$ws.Range("D16").Value = '4.278.04'
$ws.Range("E16").Value = '  -2.75%  '

# Row 17
$ws.Range("D17").Value = '95.257.85'
$ws.Range("E17").Value = '  -1.78%  '

# Row 18
$ws.Range("D18").Value = '''23.18'
$ws.Range("E18").Value = '  +24.05%  '

# Row 19
$ws.Range("D19").Value = '''8.91'
$ws.Range("E19").Value = '  -2.37%  '

# Row 20
$ws.Range("D20").Value = '''13.83'
$ws.Range("E20").Value = '  +6.65%  '

# Row 21
$ws.Range("D21").Value = '3.601.28'
$ws.Range("E21").Value = '  -3.02%  '

# Row 22
$ws.Range("D22").Value = '''0.292'
$ws.Range("E22").Value = '  +48.24%  '

# Row 23
$ws.Range("D23").Value = '''0.544'
$ws.Range("E23").Value = '  +3.42%  '

# Row 24
$ws.Range("D24").Value = '''514.72'
$ws.Range("E24").Value = '  -1.68%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '''127.36'
$ws.Range("E25").Value = '  +18.45%  '

# Row 26
$ws.Range("B26").Value = 'SuiNetwork'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D26").Value = '''3.24'
$ws.Range("E26").Value = '  -6.18%  '

# Row 27
$ws.Range("D27").Value = '''0.0000202'
$ws.Range("E27").Value = '  -9.74%  '

# Row 28
$ws.Range("D28").Value = '''6.80'
$ws.Range("E28").Value = '  -1.22%  '

# Row 29
$ws.Range("D29").Value = '3.772.34'
$ws.Range("E29").Value = '  -3.38%  '

# Row 30
$ws.Range("D30").Value = '''12.79'
$ws.Range("E30").Value = '  -5.21%  '

# Row 31
$ws.Range("D31").Value = '''13.06'
$ws.Range("E31").Value = '  +3.76%  '

# Row 32
$ws.Range("D32").Value = '''3.07'
$ws.Range("E32").Value = '  +1.41%  '

# Row 33
$ws.Range("E33").Value = '  +0.10%  '

# Row 34
$ws.Range("D34").Value = '''0.623'
$ws.Range("E34").Value = '  +5.08%  '

# Row 35
$ws.Range("D35").Value = '''0.181'
$ws.Range("E35").Value = '  -4.79%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''32.79'
$ws.Range("E36").Value = '  +0.72%  '

# Row 37
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  +0.10%  '

# Row 38
$ws.Range("D38").Value = '''1.77'
$ws.Range("E38").Value = '  -3.24%  '

# Row 39
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").Value = '''0.544'
$ws.Range("E39").Value = '  +7.82%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  -0.01%  '

# Row 41
$ws.Range("D41").Value = '''7.18'
$ws.Range("E41").Value = '  +5.39%  '

# Row 42
$ws.Range("D42").Value = '''8.34'
$ws.Range("E42").Value = '  -4.44%  '

# Row 43
$ws.Range("D43").Value = '''583.88'
$ws.Range("E43").Value = '  -9.14%  '

# Row 44
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '''0.0522'
$ws.Range("E44").Value = '  +14.41%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''41.68'
$ws.Range("E45").Value = '  +4.04%  '

# Row 46
$ws.Range("D46").Value = '''0.156'
$ws.Range("E46").Value = '  -6.09%  '

# Row 47
$ws.Range("D47").Value = '''0.962'
$ws.Range("E47").Value = '  +0.19%  '

# Row 48
$ws.Range("D48").Value = '''1.94'
$ws.Range("E48").Value = '  -4.25%  '

# Row 49
$ws.Range("D49").Value = '''9.19'
$ws.Range("E49").Value = '  +5.38%  '

# Row 50
$ws.Range("D50").Value = '''234.04'
$ws.Range("E50").Value = '  +13.66%  '

# Row 51
$ws.Range("D51").Value = '''23.50'
$ws.Range("E51").Value = '  -0.44%  '
